$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Row 6: update labels on D6/H6/L6 and clear E6
$ws.Range("D6").Value = "Cost increaes (kEUROS)"
$ws.Range("H6").Value = "Cost increaes (kEUROS)"
$ws.Range("L6").Value = "Cost increaes (kEUROS)"
$ws.Range("E6").ClearContents()

# Rows 7-15: move cost formula from E/I/M into D/H/L, then clear E/I/M
for ($r = 7; $r -le 15; $r++) {
    $ws.Range("D$r").Formula = "=C$r*rfonfcost_per_station_per_km/1000"
    $ws.Range("H$r").Formula = "=G$r*rfonfcost_per_station_per_km/1000"
    $ws.Range("L$r").Formula = "=K$r*rfonfcost_per_station_per_km/1000"
    $ws.Range("E$r").ClearContents()
    $ws.Range("I$r").ClearContents()
    $ws.Range("M$r").ClearContents()
}

# Row 41: remove spurious MST total row contents (keep B41/C41 cells with style, cleared)
$ws.Range("A41").ClearContents()
$ws.Range("B41").ClearContents()
$ws.Range("C41").ClearContents()
$ws.Range("D41").ClearContents()

$ws.Range("L38").Select()

Write-Host "done"
